# "Add MS Band Code" -- insert a new "Adaptive Code" slide right after the
# existing "Adaptive UI - ???????" placeholder slide (slide 6), using the
# "Title and Content" layout, and fill in its title/body text.

$p = $ppt.ActivePresentation

# Slide 6 currently holds the "Adaptive" / "UI - ???????" placeholder slide.
# Insert the new slide immediately after it (new position 7), using
# ppLayoutText (2) == "Title and Content" layout.
$new = $p.Slides.Add(7, 2)

# Title placeholder
$new.Shapes.Item(1).TextFrame.TextRange.Text = "Adaptive Code"

# Body / content placeholder
$body = $new.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Most of the Windows Runtime APIs your Universal 8.1 app already calls are implemented in the set of APIs known as the universal device family. But, some are implemented in extension SDKs, and Visual Studio only recognizes APIs that are implemented by your app's target device family or by any extension SDKs that you have referenced."
$body.ParagraphFormat.Bullet.Visible = $false
